# Updates Sheet1 "Price" (D) / "Volume(1h)" (E) columns to match the refreshed
# cryptos snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Several "Price" cells are numeric-looking strings (e.g. "5.94").
    # Plain `.Value = "5.94"` gets auto-coerced to a real number by Excel,
    # which does not match the source data (plain text cells, no thousands
    # grouping). Forcing NumberFormat to Text keeps the literal string, and
    # ClearFormats() afterwards drops the number-format stamp again so the
    # cells style index is left exactly as it was.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = '64.081.87'
$ws.Cells.Item(2, 5).Value = '  +0.43%  '
$ws.Cells.Item(3, 4).Value = '3.141.66'
$ws.Cells.Item(3, 5).Value = '  +0.81%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
Set-TextValue $ws.Cells.Item(5, 4) '590.36'
$ws.Cells.Item(5, 5).Value = '  +0.71%  '
Set-TextValue $ws.Cells.Item(6, 4) '145.56'
$ws.Cells.Item(6, 5).Value = '  -0.50%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '3.132.88'
$ws.Cells.Item(8, 5).Value = '  +0.74%  '
$ws.Cells.Item(9, 5).Value = '  -0.22%  '
$ws.Cells.Item(10, 5).Value = '  -0.15%  '
Set-TextValue $ws.Cells.Item(11, 4) '5.94'
$ws.Cells.Item(11, 5).Value = '  +2.82%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.458'
$ws.Cells.Item(12, 5).Value = '  -1.40%  '
Set-TextValue $ws.Cells.Item(13, 4) '0.0000247'
$ws.Cells.Item(13, 5).Value = '  -1.53%  '
Set-TextValue $ws.Cells.Item(14, 4) '37.45'
$ws.Cells.Item(14, 5).Value = '  +0.64%  '
$ws.Cells.Item(15, 4).Value = '3.656.89'
$ws.Cells.Item(15, 5).Value = '  +0.67%  '
$ws.Cells.Item(16, 5).Value = '  -1.21%  '
Set-TextValue $ws.Cells.Item(17, 4) '7.31'
$ws.Cells.Item(17, 5).Value = '  +2.43%  '
$ws.Cells.Item(18, 4).Value = '63.892.37'
$ws.Cells.Item(18, 5).Value = '  +0.29%  '
$ws.Cells.Item(19, 4).Value = '3.136.75'
$ws.Cells.Item(19, 5).Value = '  +0.79%  '
Set-TextValue $ws.Cells.Item(20, 4) '468.00'
$ws.Cells.Item(20, 5).Value = '  +0.73%  '
Set-TextValue $ws.Cells.Item(21, 4) '14.35'
$ws.Cells.Item(21, 5).Value = '  +0.29%  '
Set-TextValue $ws.Cells.Item(22, 4) '0.733'
Set-TextValue $ws.Cells.Item(23, 4) '7.56'
$ws.Cells.Item(23, 5).Value = '  +0.53%  '
Set-TextValue $ws.Cells.Item(24, 4) '2.34'
$ws.Cells.Item(24, 5).Value = '  +8.26%  '
Set-TextValue $ws.Cells.Item(25, 4) '12.99'
$ws.Cells.Item(25, 5).Value = '  -1.30%  '
Set-TextValue $ws.Cells.Item(26, 4) '81.46'
$ws.Cells.Item(26, 5).Value = '  -0.43%  '
$ws.Cells.Item(27, 5).Value = '  +0.10%  '
Set-TextValue $ws.Cells.Item(28, 4) '9.97'
$ws.Cells.Item(28, 5).Value = '  +11.79%  '
Set-TextValue $ws.Cells.Item(29, 4) '7.49'
$ws.Cells.Item(29, 5).Value = '  +8.98%  '
$ws.Cells.Item(30, 5).Value = '  +0.66%  '
$ws.Cells.Item(31, 5).Value = '  +0.32%  '
$ws.Cells.Item(32, 5).Value = '  +0.14%  '
Set-TextValue $ws.Cells.Item(33, 4) '27.61'
$ws.Cells.Item(33, 5).Value = '  +2.24%  '
$ws.Cells.Item(34, 5).Value = '  +1.35%  '
$ws.Cells.Item(35, 4).Value = '0.0₃0849'
$ws.Cells.Item(35, 5).Value = '  -2.48%  '
$ws.Cells.Item(36, 5).Value = '  +1.25%  '
Set-TextValue $ws.Cells.Item(37, 4) '6.16'
$ws.Cells.Item(37, 5).Value = '  +1.57%  '
Set-TextValue $ws.Cells.Item(38, 4) '2.31'
$ws.Cells.Item(38, 5).Value = '  -2.20%  '
Set-TextValue $ws.Cells.Item(39, 4) '3.21'
$ws.Cells.Item(39, 5).Value = '  -5.98%  '
Set-TextValue $ws.Cells.Item(40, 4) '51.39'
$ws.Cells.Item(40, 5).Value = '  +0.89%  '
Set-TextValue $ws.Cells.Item(41, 4) '9.31'
$ws.Cells.Item(41, 5).Value = '  +7.26%  '
Set-TextValue $ws.Cells.Item(42, 4) '456.11'
$ws.Cells.Item(42, 5).Value = '  +1.72%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.294'
$ws.Cells.Item(43, 5).Value = '  +6.10%  '
Set-TextValue $ws.Cells.Item(44, 4) '0.0373'
$ws.Cells.Item(44, 5).Value = '  +0.39%  '
$ws.Cells.Item(45, 4).Value = '2.918.69'
$ws.Cells.Item(45, 5).Value = '  +1.43%  '
Set-TextValue $ws.Cells.Item(46, 4) '40.21'
$ws.Cells.Item(46, 5).Value = '  +12.48%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.108'
$ws.Cells.Item(47, 5).Value = '  -2.77%  '
Set-TextValue $ws.Cells.Item(48, 4) '133.78'
$ws.Cells.Item(48, 5).Value = '  +8.16%  '
$ws.Cells.Item(49, 5).Value = '  +0.00%  '
$ws.Cells.Item(50, 5).Value = '  +2.68%  '
$ws.Cells.Item(51, 5).Value = '  -0.47%  '
